$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reaction-term formulas (column E) referenced a parameter called
# "beta"; it has been renamed to "G". Update the two text cells that spell
# the expression out (they are plain text, not real spreadsheet formulas).
$ws.Range("E2").Value = "-b1*R+((K)/(1+(G*T)))"
$ws.Range("E3").Value = "g1*R-b2*L+((K)/(1+(G*T)))"

# Rename the parameter itself in the "Parameters" column.
$ws.Range("D8").Value = "G"

# Move the sheet's active cell/selection to E3.
$ws.Range("E3").Select() | Out-Null
